$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gaz")

# Set the C3 value (Mars/Avr 2017 gas consumption)
$ws.Range("C3").Value = 23.45

# Move the active cell / selection to C10
$ws.Range("C10").Select()
